# Ontwerp Summatieve opdracht 2.docx - apply commit "Add files via upload"
#
# Change 1: the document title paragraph switches from heading style
#           "Kop1" to "Kop2".
# Change 2: two new paragraphs are added right before the very last
#           (empty) paragraph of the document:
#             - a "Kop2" heading "Aanpasbaarheid"
#             - a body paragraph describing that both FSM's are easy to
#               adapt, with the usual Dutch-spellcheck <w:proofErr/>
#               bookmarks around the words Word's proofer would flag.

$d = $word.ActiveDocument

# --- Change 1 : Kop1 -> Kop2 on the title paragraph -------------------
$title = $d.Paragraphs.First
$title.Style = "Kop2"

# --- Change 2 : insert the new "Aanpasbaarheid" section ----------------
# Build the two paragraphs as literal WordprocessingML and splice them
# in front of the last (empty) paragraph via Range.InsertXML, so the
# run/proofErr layout matches exactly what Word itself would have
# produced while typing this text.
$newParagraphsXml = '<w:p><w:pPr><w:pStyle w:val="Kop2"/></w:pPr><w:r><w:t>Aanpasbaarheid</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t xml:space="preserve">Allebei de </w:t></w:r>' + `
'<w:proofErr w:type="spellStart"/><w:r><w:t>FSM' + [char]0x2019 + 's</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
'<w:r><w:t xml:space="preserve"> zijn makkelijk aan te passen in </w:t></w:r>' + `
'<w:proofErr w:type="spellStart"/><w:r><w:t>CreateDobbelFSM</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
'<w:r><w:t xml:space="preserve"> en </w:t></w:r>' + `
'<w:proofErr w:type="spellStart"/><w:r><w:t>CreateTekstFSM</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
'<w:r><w:t xml:space="preserve">. </w:t></w:r></w:p>'

$package = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body>' + $newParagraphsXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$lastPara = $d.Paragraphs.Last
$null = $lastPara.Range.InsertXML($package)
